$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (BC1) - 24/05/20 (one day after BB1's 23/05/20)
$ws.Range("BC1").Value = 43975

# New case counts for column BC (rows 2-19)
$ws.Range("BC2").Value = 109
$ws.Range("BC3").Value = 86
$ws.Range("BC4").Value = 42
$ws.Range("BC5").Value = 35
$ws.Range("BC6").Value = 5
$ws.Range("BC7").Value = 2543
$ws.Range("BC8").Value = 23
$ws.Range("BC9").Value = 771
$ws.Range("BC10").Value = 0
$ws.Range("BC11").Value = 13
$ws.Range("BC12").Value = 1
$ws.Range("BC13").Value = 9
$ws.Range("BC14").Value = 12
$ws.Range("BC15").Value = 3
$ws.Range("BC16").Value = 11
$ws.Range("BC17").Value = 35
$ws.Range("BC18").Value = 45
$ws.Range("BC19").Value = 207

# Totals row
$ws.Range("BC20").Formula = "=SUM(BC2:BC19)"

# Move / update the active selection to match the post-edit state
$ws.Range("BD14").Select()
